$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Definition" sheet: the list of node names (rows 7-18, column A) was
#    reshuffled into a new order. "Water" (row 8) stays put; the rest move.
# ---------------------------------------------------------------------------
$wsDef = $wb.Worksheets.Item("Definition")

$wsDef.Range("A7").Value = "District_Heating"
$wsDef.Range("A9").Value = "Vaporized_Carbon_Dioxide"
$wsDef.Range("A10").Value = "E-Methanol_Kasso"
$wsDef.Range("A11").Value = "Waste_Heat"
$wsDef.Range("A12").Value = "E-Methanol_storage_Kasso"
$wsDef.Range("A13").Value = "Power_Wholesale"
$wsDef.Range("A14").Value = "Carbon_Dioxide"
$wsDef.Range("A15").Value = "Hydrogen_Kasso"
$wsDef.Range("A16").Value = "Hydrogen_storage_Kasso"
$wsDef.Range("A17").Value = "Raw_Methanol"
$wsDef.Range("A18").Value = "Power_Kasso"

# ---------------------------------------------------------------------------
# 2) "Nodes" sheet: same reshuffle, but the whole row of attributes
#    (balance_type, has_state, node_state_cap, frac_state_loss,
#    node_slack_penalty) travels together with each node's name. Only the
#    cells whose effective value actually changes are touched.
# ---------------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("Nodes")

$wsNodes.Range("A2").Value = "District_Heating"
$wsNodes.Range("C2").Value = "balance_type_none"

$wsNodes.Range("A4").Value = "Vaporized_Carbon_Dioxide"

$wsNodes.Range("A5").Value = "E-Methanol_Kasso"
$wsNodes.Range("C5").Value = "balance_type_node"
$wsNodes.Range("G5").Value = 100000

$wsNodes.Range("A6").Value = "Waste_Heat"
$wsNodes.Range("G6").ClearContents()

$wsNodes.Range("A7").Value = "E-Methanol_storage_Kasso"
$wsNodes.Range("C7").Value = "balance_type_node"
$wsNodes.Range("D7").Value = $true
$wsNodes.Range("E7").Value = 100000
$wsNodes.Range("F7").Value = 0
$wsNodes.Range("G7").Value = 100000

$wsNodes.Range("A8").Value = "Power_Wholesale"
$wsNodes.Range("C8").Value = "balance_type_none"
$wsNodes.Range("D8").ClearContents()
$wsNodes.Range("E8").ClearContents()
$wsNodes.Range("F8").ClearContents()
$wsNodes.Range("G8").ClearContents()

$wsNodes.Range("A9").Value = "Carbon_Dioxide"

$wsNodes.Range("A10").Value = "Hydrogen_Kasso"

$wsNodes.Range("A11").Value = "Hydrogen_storage_Kasso"

$wsNodes.Range("A12").Value = "Raw_Methanol"

$wsNodes.Range("A13").Value = "Power_Kasso"

# ---------------------------------------------------------------------------
# 3) "Object__to_from_node" sheet: new fom_cost / vom_cost columns (G, H).
# ---------------------------------------------------------------------------
$wsObj = $wb.Worksheets.Item("Object__to_from_node")

$wsObj.Range("G1").Value = "fom_cost"
$wsObj.Range("H1").Value = "vom_cost"

# Solar_Plant_Kasso unit__to_node Power_Kasso -> fom_cost 100
$wsObj.Range("G2").Value = 100
# Electrolyzer unit__to_node Hydrogen_Kasso -> vom_cost 1
$wsObj.Range("H4").Value = 1
# power_line_Wholesale_Kasso connection__to_node Power_Kasso -> fom_cost 100
$wsObj.Range("G16").Value = 100
# pipeline_storage_hydrogen connection__to_node Hydrogen_storage_Kasso -> vom_cost 1
$wsObj.Range("H20").Value = 1

# ---------------------------------------------------------------------------
# 4) "Variable_Eff" sheet: header rename, a value correction, and removal
#    of two now-unused trailing rows.
# ---------------------------------------------------------------------------
$wsVar = $wb.Worksheets.Item("Variable_Eff")

$wsVar.Range("C1").Value = "unit__to_node"
$wsVar.Range("C3").Value = "Hydrogen_Kasso"
$wsVar.Range("B6").Value = 1

$wsVar.Rows.Item(9).Delete() | Out-Null
$wsVar.Rows.Item(9).Delete() | Out-Null
